# Applies updated market-price / profit figures to the Leve profit sheets
# (per scheduled market-data refresh). Values below are taken verbatim
# from the refreshed dataset for each (sheet, row) pair.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1543.909
$ws.Range("I98").Value = 1198.45
$ws.Range("K98").Value = 1198.45
$ws.Range("M98").Value = 299.55
# Row 122
$ws.Range("H122").Value = 1543.909
$ws.Range("I122").Value = 1198.45
$ws.Range("K122").Value = 3595.35
$ws.Range("M122").Value = -1145.35
# Row 137
$ws.Range("H137").Value = 6533.8613
$ws.Range("I137").Value = 7385.522
$ws.Range("J137").Value = 5027.077
$ws.Range("K137").Value = 22156.566
$ws.Range("L137").Value = 15081.231
$ws.Range("M137").Value = -19606.566
$ws.Range("N137").Value = -20181.231
# Row 138
$ws.Range("H138").Value = 28572692
$ws.Range("I138").Value = 31250956
$ws.Range("K138").Value = 93752868
$ws.Range("M138").Value = -93747728

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1417.87
$ws.Range("I32").Value = 1401.899
$ws.Range("K32").Value = 1401.899
$ws.Range("M32").Value = -1114.899
# Row 61
$ws.Range("H61").Value = 98761.664
$ws.Range("I61").Value = 7312.6
$ws.Range("K61").Value = 7312.6
$ws.Range("M61").Value = -7100.6
# Row 63
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1314
$ws.Range("N63").Value = -3372
# Row 66
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -6568
$ws.Range("N66").Value = -16864
# Row 74
$ws.Range("H74").Value = 1552
$ws.Range("I74").Value = 1533.5
$ws.Range("K74").Value = 1533.5
$ws.Range("M74").Value = -659.5
# Row 77
$ws.Range("H77").Value = 1552
$ws.Range("I77").Value = 1533.5
$ws.Range("K77").Value = 7667.5
$ws.Range("M77").Value = -3299.5
# Row 96
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
# Row 104
$ws.Range("H104").Value = 29933.334
$ws.Range("J104").Value = 29933.334
$ws.Range("L104").Value = 29933.334
$ws.Range("N104").Value = -36921.334
# Row 132
$ws.Range("H132").Value = 126145.25
$ws.Range("I132").Value = 67887.336
$ws.Range("J132").Value = 1000014
$ws.Range("K132").Value = 203662.008
$ws.Range("L132").Value = 3000042
$ws.Range("M132").Value = -201132.008
$ws.Range("N132").Value = -3005102
# Row 135
$ws.Range("H135").Value = 78499.25
$ws.Range("J135").Value = 78499.25
$ws.Range("L135").Value = 78499.25
$ws.Range("N135").Value = -88639.25
# Row 136
$ws.Range("H136").Value = 98761.664
$ws.Range("I136").Value = 7312.6
$ws.Range("K136").Value = 21937.8
$ws.Range("M136").Value = -19387.8
# Row 141
$ws.Range("H141").Value = 250664.67
$ws.Range("J141").Value = 250664.67
$ws.Range("L141").Value = 250664.67
$ws.Range("N141").Value = -261024.67

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Range("H76").Value = 52750
$ws.Range("J76").Value = 52750
$ws.Range("L76").Value = 52750
$ws.Range("N76").Value = -53380
# Row 79
$ws.Range("H79").Value = 52750
$ws.Range("J79").Value = 52750
$ws.Range("L79").Value = 52750
$ws.Range("N79").Value = -54934
# Row 95
$ws.Range("H95").Value = 99999.5
$ws.Range("J95").Value = 99999.5
$ws.Range("L95").Value = 99999.5
$ws.Range("N95").Value = -105491.5
# Row 134
$ws.Range("H134").Value = 7134.0713
$ws.Range("I134").Value = 3806.5908
$ws.Range("K134").Value = 11419.7724
$ws.Range("M134").Value = -8884.7724

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1896.7715
$ws.Range("I122").Value = 1569.64
$ws.Range("J122").Value = 2714.6
$ws.Range("K122").Value = 4708.92
$ws.Range("L122").Value = 8143.799999999999
$ws.Range("M122").Value = -2258.92
$ws.Range("N122").Value = -13043.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1350.1154
$ws.Range("I5").Value = 1326.7
$ws.Range("J5").Value = 1364.75
$ws.Range("K5").Value = 3980.1
$ws.Range("L5").Value = 4094.25
$ws.Range("M5").Value = -3868.1
$ws.Range("N5").Value = -4318.25
# Row 68
$ws.Range("H68").Value = 883.4167
$ws.Range("I68").Value = 993.625
$ws.Range("J68").Value = 663
$ws.Range("K68").Value = 2980.875
$ws.Range("L68").Value = 1989
$ws.Range("M68").Value = -2169.875
$ws.Range("N68").Value = -3611
# Row 71
$ws.Range("H71").Value = 883.4167
$ws.Range("I71").Value = 993.625
$ws.Range("J71").Value = 663
$ws.Range("K71").Value = 8942.625
$ws.Range("L71").Value = 5967
$ws.Range("M71").Value = -4886.625
$ws.Range("N71").Value = -14079
# Row 103
$ws.Range("H103").Value = 96.666664
$ws.Range("I103").Value = 96.666664
$ws.Range("K103").Value = 289.999992
$ws.Range("M103").Value = 589.000008
# Row 135
$ws.Range("H135").Value = 1350.1154
$ws.Range("I135").Value = 1326.7
$ws.Range("J135").Value = 1364.75
$ws.Range("K135").Value = 11940.3
$ws.Range("L135").Value = 12282.75
$ws.Range("M135").Value = -9405.300000000001
$ws.Range("N135").Value = -17352.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 31397.086
$ws.Range("I132").Value = 2039.6
$ws.Range("J132").Value = 104790.8
$ws.Range("K132").Value = 6118.799999999999
$ws.Range("L132").Value = 314372.4
$ws.Range("M132").Value = -3588.799999999999
$ws.Range("N132").Value = -319432.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 130
$ws.Range("H130").Value = 147996.25
$ws.Range("J130").Value = 147996.25
$ws.Range("L130").Value = 147996.25
$ws.Range("N130").Value = -158036.25
# Row 132
$ws.Range("H132").Value = 31077.092
$ws.Range("I132").Value = 33381.41
$ws.Range("J132").Value = 1697
$ws.Range("K132").Value = 100144.23
$ws.Range("L132").Value = 5091
$ws.Range("M132").Value = -97614.23000000001
$ws.Range("N132").Value = -10151
# Row 136
$ws.Range("H136").Value = 8317.105
$ws.Range("I136").Value = 9743.272000000001
$ws.Range("K136").Value = 29229.816
$ws.Range("M136").Value = -26679.816

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5147.4165
$ws.Range("I132").Value = 4099.706
$ws.Range("J132").Value = 7691.857
$ws.Range("K132").Value = 12299.118
$ws.Range("L132").Value = 23075.571
$ws.Range("M132").Value = -9769.118
$ws.Range("N132").Value = -28135.571
# Row 136
$ws.Range("H136").Value = 2419.6667
$ws.Range("I136").Value = 2270.6553
$ws.Range("J136").Value = 2851.8
$ws.Range("K136").Value = 6811.965899999999
$ws.Range("L136").Value = 8555.400000000001
$ws.Range("M136").Value = -4261.965899999999
$ws.Range("N136").Value = -13655.4

Write-Output "Applied Leve profit updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"